$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1890.5555
$ws.Range("J17").Value = 1925.2941
$ws.Range("L17").Value = 5775.8823
$ws.Range("N17").Value = -6111.8823
# Row 29
$ws.Range("H29").Value = 4308.3335
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 4473.75
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 13421.25
$ws.Range("M29").Value = -2719
$ws.Range("N29").Value = -13983.25
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 111
$ws.Range("H111").Value = 88386.086
$ws.Range("J111").Value = 7326
$ws.Range("L111").Value = 21978
$ws.Range("N111").Value = -28112
# Row 113
$ws.Range("H113").Value = 5922.7144
$ws.Range("I113").Value = 5268.778
$ws.Range("J113").Value = 7099.8
$ws.Range("K113").Value = 5268.778
$ws.Range("L113").Value = 7099.8
$ws.Range("M113").Value = -2014.778
$ws.Range("N113").Value = -13607.8
# Row 132
$ws.Range("H132").Value = 2743.074
$ws.Range("I132").Value = 2824.158
$ws.Range("J132").Value = 2550.5
$ws.Range("K132").Value = 8472.474
$ws.Range("L132").Value = 7651.5
$ws.Range("M132").Value = -5942.474
$ws.Range("N132").Value = -12711.5
# Row 137
$ws.Range("H137").Value = 694937.5600000001
$ws.Range("I137").Value = 528570.4
$ws.Range("K137").Value = 1585711.2
$ws.Range("M137").Value = -1583161.2

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3908.0833
$ws.Range("I61").Value = 2499.625
$ws.Range("K61").Value = 2499.625
$ws.Range("M61").Value = -2287.625
# Row 122
$ws.Range("H122").Value = 3308.2974
$ws.Range("I122").Value = 1816.6086
$ws.Range("K122").Value = 5449.825800000001
$ws.Range("M122").Value = -2999.825800000001
# Row 132
$ws.Range("H132").Value = 2698.8
$ws.Range("I132").Value = 2200.963
$ws.Range("K132").Value = 6602.889000000001
$ws.Range("M132").Value = -4072.889000000001
# Row 136
$ws.Range("H136").Value = 3908.0833
$ws.Range("I136").Value = 2499.625
$ws.Range("K136").Value = 7498.875
$ws.Range("M136").Value = -4948.875

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 4018.6667
$ws.Range("I99").Value = 3318.125
$ws.Range("K99").Value = 3318.125
$ws.Range("M99").Value = -1820.125

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 35850.25
$ws.Range("I31").Value = 1538.0588
$ws.Range("K31").Value = 1538.0588
$ws.Range("M31").Value = -1243.0588
# Row 34
$ws.Range("H34").Value = 35850.25
$ws.Range("I34").Value = 1538.0588
$ws.Range("K34").Value = 1538.0588
$ws.Range("M34").Value = -1336.0588
# Row 58
$ws.Range("H58").Value = 375919.97
$ws.Range("I58").Value = 560352.2
$ws.Range("J58").Value = 7055.5557
$ws.Range("K58").Value = 560352.2
$ws.Range("L58").Value = 7055.5557
$ws.Range("M58").Value = -560149.2
$ws.Range("N58").Value = -7461.5557
# Row 99
$ws.Range("H99").Value = 4775.9614
$ws.Range("J99").Value = 5490.8335
$ws.Range("L99").Value = 5490.8335
$ws.Range("N99").Value = -8486.833500000001
# Row 122
$ws.Range("H122").Value = 3419.1738
$ws.Range("I122").Value = 1991.6923
$ws.Range("J122").Value = 5274.9
$ws.Range("K122").Value = 5975.0769
$ws.Range("L122").Value = 15824.7
$ws.Range("M122").Value = -3525.0769
$ws.Range("N122").Value = -20724.7
# Row 126
$ws.Range("H126").Value = 4775.9614
$ws.Range("J126").Value = 5490.8335
$ws.Range("L126").Value = 16472.5005
$ws.Range("N126").Value = -21412.5005
# Row 132
$ws.Range("H132").Value = 4166.6943
$ws.Range("I132").Value = 3849.5386
$ws.Range("J132").Value = 4991.3
$ws.Range("K132").Value = 11548.6158
$ws.Range("L132").Value = 14973.9
$ws.Range("M132").Value = -9018.6158
$ws.Range("N132").Value = -20033.9
# Row 136
$ws.Range("H136").Value = 375919.97
$ws.Range("I136").Value = 560352.2
$ws.Range("J136").Value = 7055.5557
$ws.Range("K136").Value = 1681056.6
$ws.Range("L136").Value = 21166.6671
$ws.Range("M136").Value = -1678506.6
$ws.Range("N136").Value = -26266.6671

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
# Row 13
$ws.Range("H13").Value = 2000
$ws.Range("I13").Value = 1466.5
$ws.Range("J13").Value = 2640.2
$ws.Range("K13").Value = 4399.5
$ws.Range("L13").Value = 7920.599999999999
$ws.Range("M13").Value = -4231.5
$ws.Range("N13").Value = -8256.599999999999
# Row 17
$ws.Range("H17").Value = 1025
$ws.Range("I17").Value = 300
$ws.Range("K17").Value = 900
$ws.Range("M17").Value = -731
# Row 55
$ws.Range("H55").Value = 42989.684
$ws.Range("J55").Value = 49993.75
$ws.Range("L55").Value = 149981.25
$ws.Range("N55").Value = -150335.25
# Row 132
$ws.Range("H132").Value = 735680.5600000001
$ws.Range("J132").Value = 1670242.5
$ws.Range("L132").Value = 15032182.5
$ws.Range("N132").Value = -15037242.5
# Row 137
$ws.Range("H137").Value = 2115.8635
$ws.Range("J137").Value = 5199.8
$ws.Range("L137").Value = 15599.4
$ws.Range("N137").Value = -25799.4
# Row 139
$ws.Range("H139").Value = 4699.5
$ws.Range("I139").Value = 4699.5
$ws.Range("K139").Value = 14098.5
$ws.Range("M139").Value = -8958.5

$ws = $wb.Worksheets.Item("GSM")
# Row 17
$ws.Range("H17").Value = 8267.200000000001
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 10309
$ws.Range("K17").Value = 100
$ws.Range("L17").Value = 10309
$ws.Range("M17").Value = 68
$ws.Range("N17").Value = -10645
# Row 23
$ws.Range("H23").Value = 989
$ws.Range("J23").Value = 989
$ws.Range("L23").Value = 989
$ws.Range("N23").Value = -1435
# Row 80
$ws.Range("H80").Value = 2504051
$ws.Range("J80").Value = 3335467
$ws.Range("L80").Value = 3335467
$ws.Range("N80").Value = -3337463
# Row 83
$ws.Range("H83").Value = 2504051
$ws.Range("J83").Value = 3335467
$ws.Range("L83").Value = 16677335
$ws.Range("N83").Value = -16687319
# Row 132
$ws.Range("H132").Value = 359724.4
$ws.Range("I132").Value = 561404.7
$ws.Range("J132").Value = 80474.84
$ws.Range("K132").Value = 1684214.1
$ws.Range("L132").Value = 241424.52
$ws.Range("M132").Value = -1681684.1
$ws.Range("N132").Value = -246484.52

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 349228.62
$ws.Range("I7").Value = 4450.533
$ws.Range("J7").Value = 718633.7
$ws.Range("K7").Value = 4450.533
$ws.Range("L7").Value = 718633.7
$ws.Range("M7").Value = -4338.533
$ws.Range("N7").Value = -718857.7
# Row 40
$ws.Range("H40").Value = 34552.195
$ws.Range("I40").Value = 44930.176
$ws.Range("K40").Value = 44930.176
$ws.Range("M40").Value = -44794.176
# Row 126
$ws.Range("H126").Value = 349228.62
$ws.Range("I126").Value = 4450.533
$ws.Range("J126").Value = 718633.7
$ws.Range("K126").Value = 13351.599
$ws.Range("L126").Value = 2155901.1
$ws.Range("M126").Value = -10881.599
$ws.Range("N126").Value = -2160841.1
# Row 132
$ws.Range("H132").Value = 4669.5713
$ws.Range("I132").Value = 2249.75
$ws.Range("J132").Value = 5637.5
$ws.Range("K132").Value = 6749.25
$ws.Range("L132").Value = 16912.5
$ws.Range("M132").Value = -4219.25
$ws.Range("N132").Value = -21972.5
# Row 136
$ws.Range("H136").Value = 1182134.4
$ws.Range("I136").Value = 1670815.1
$ws.Range("K136").Value = 5012445.300000001
$ws.Range("M136").Value = -5009895.300000001

$ws = $wb.Worksheets.Item("WVR")
# Row 93
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992
# Row 126
$ws.Range("H126").Value = 1629.2858
$ws.Range("I126").Value = 1367.7307
$ws.Range("J126").Value = 2384.889
$ws.Range("K126").Value = 4103.1921
$ws.Range("L126").Value = 7154.667
$ws.Range("M126").Value = -1633.1921
$ws.Range("N126").Value = -12094.667
# Row 132
$ws.Range("H132").Value = 38532.07
$ws.Range("I132").Value = 3470.25
$ws.Range("J132").Value = 116447.22
$ws.Range("K132").Value = 10410.75
$ws.Range("L132").Value = 349341.66
$ws.Range("M132").Value = -7880.75
$ws.Range("N132").Value = -354401.66
# Row 136
$ws.Range("H136").Value = 10170793
$ws.Range("I136").Value = 17188410
$ws.Range("K136").Value = 51565230
$ws.Range("M136").Value = -51562680
